$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column S (19th column), pushing the existing
# "Hoja" column (and its per-row "POAI_2025" values) from S to T.
$ws.Columns("S").Insert()

# New column header for the inserted column S.
$ws.Range("S1").Value = "Col19"

# Enlace Técnico for row 6 (previously blank).
$ws.Range("M6").Value = "ANDREA GONZALEZ"

# Responsable corrections.
$ws.Range("L8").Value = "SARA DIANA URBANO"
$ws.Range("L11").Value = "LUZ MIRYAN Y WILLIAN MOSQUERA"
